# 13项目计划表.xlsx — add two new weekly-plan blocks (rows 114-122 and 125-133)
# mirroring the existing block layout, fill in C106:C110 completion values, and
# update the current selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the previously-empty "completion" values for the last existing
#    block (rows 106-110).
# ---------------------------------------------------------------------------
$ws.Range("C106").Value = 1
$ws.Range("C107").Value = 0.5
$ws.Range("C108").Value = 1
$ws.Range("C109").Value = 0.5
$ws.Range("C110").Value = 1

# ---------------------------------------------------------------------------
# 2) Clone the formatting of the last full block (rows 104-112) onto the two
#    new blocks so every cell picks up the same fonts/borders/alignment/
#    number-formats already used throughout the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A104:D112").Copy()
$ws.Range("A114:D122").PasteSpecial(-4122)

$ws.Range("A104:D112").Copy()
$ws.Range("A125:D133").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Block 1: 日期：2018.11.01 第九周周四  (rows 114-122)
# ---------------------------------------------------------------------------
$ws.Range("A114").Value = "日期：2018.11.01 第九周周四"

$ws.Range("A115").Value = "组员"
$ws.Range("B115").Value = "计划内容"
$ws.Range("C115").Value = "完成情况"
$ws.Range("D115").Value = "备注"

$ws.Range("A116").Value = "王伟锋"
$ws.Range("B116").Value = "重构后台框架为spring mvc"
$ws.Range("C116").Value = 1

$ws.Range("A117").Value = "陈升云"
$ws.Range("B117").Value = "完成个人资料的查看，修改等"
$ws.Range("C117").Value = 0.7

$ws.Range("A118").Value = "林玮成"
$ws.Range("B118").Value = "辅助app开发"
$ws.Range("C118").Value = 1

$ws.Range("A119").Value = "吴帅辰"
$ws.Range("B119").Value = "完成查看历史系统推送消息功能"
$ws.Range("C119").Value = 0.7

$ws.Range("A120").Value = "李海洋"
$ws.Range("B120").Value = "完成群组和消息界面代码的编写"
$ws.Range("C120").Value = 1

$ws.Range("A121").Value = "总结："

$ws.Range("A114:D114").Merge()
$ws.Range("A121:D122").Merge()

# ---------------------------------------------------------------------------
# 4) Block 2: 日期：2018.11.05 第十周周一  (rows 125-133)
# ---------------------------------------------------------------------------
$ws.Range("A125").Value = "日期：2018.11.05 第十周周一"

$ws.Range("A126").Value = "组员"
$ws.Range("B126").Value = "计划内容"
$ws.Range("C126").Value = "完成情况"
$ws.Range("D126").Value = "备注"

$ws.Range("A127").Value = "王伟锋"
$ws.Range("B127").Value = "完成个人头像的存储修改及个人信息的修改"

$ws.Range("A128").Value = "陈升云"
$ws.Range("B128").Value = "完成个人资料的查看，修改等"

$ws.Range("A129").Value = "林玮成"
$ws.Range("B129").Value = "辅助app开发"

$ws.Range("A130").Value = "吴帅辰"
$ws.Range("B130").Value = "完成查看历史系统推送消息功能"

$ws.Range("A131").Value = "李海洋"
$ws.Range("B131").Value = "完成获取个人所有群信息"

$ws.Range("A132").Value = "总结："

$ws.Range("A125:D125").Merge()
$ws.Range("A132:D133").Merge()

# ---------------------------------------------------------------------------
# 5) Restore the view: scroll near the bottom of the sheet and leave the
#    selection where the author left off.
# ---------------------------------------------------------------------------
$ws.Range("E132").Select()

Write-Output "done"
